$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "raj001@mail.com"
$ws.Range("D3").Value = "taj002@mail.com"
$ws.Range("D4").Value = "gani003@bhai.com"

$ws.Range("G10").Select()
